$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 7
$ws.Range("C7").Value = 8000068

# Row 18
$ws.Range("B18").Value = 10000025

# Row 19
$ws.Range("B19").Value = 10000025
$ws.Range("C19").Value = 8000025

# Row 20
$ws.Range("B20").Value = 10000025
$ws.Range("C20").Value = 8000040

# Row 21
$ws.Range("B21").Value = 10000025
$ws.Range("C21").Value = 8000043

# Row 22
$ws.Range("B22").Value = 10000025
$ws.Range("C22").Value = 8000044

# Row 23
$ws.Range("B23").Value = 10000025
$ws.Range("C23").Value = 8000047

# Row 24
$ws.Range("B24").Value = 10000025
$ws.Range("C24").Value = 8000057

# Row 25
$ws.Range("B25").Value = 10000025
$ws.Range("C25").Value = 8000063

# Row 28
$ws.Range("C28").Value = 8000069

# Row 29
$ws.Range("C29").Value = 8000075

# Row 30
$ws.Range("B30").Value = 10000034
$ws.Range("C30").Value = 8000034
$ws.Range("E30").Value = "SOREU dei Laghi"

# Row 31
$ws.Range("B31").Value = 10000034
$ws.Range("C31").Value = 8000037
$ws.Range("E31").Value = "Le SOREU"

# Row 32
$ws.Range("B32").Value = 10000038
$ws.Range("C32").Value = 8000038

# Row 33
$ws.Range("B33").Value = 10000041
$ws.Range("C33").Value = 8000041
$ws.Range("D33").Value = "Repubblica e Cantone Ticino"
$ws.Range("E33").Value = "Legge sulla protezione civile del 26 febbraio 2007"

# Row 34
$ws.Range("B34").Value = 10000041
$ws.Range("C34").Value = 8000053
$ws.Range("E34").Value = "Legge sulla protezione della popolazione (del 26 febbraio 2007)"

# Row 35
$ws.Range("B35").Value = 10000041
$ws.Range("C35").Value = 8000059
$ws.Range("E35").Value = "Regolamento sulla protezione della popolazione (RProtPop) (del 18 ottobre 2017)"

# Row 36
$ws.Range("B36").Value = 10000041
$ws.Range("C36").Value = 8000066
$ws.Range("E36").Value = "Servizio della protezione della popolazione"

# Row 37
$ws.Range("B37").Value = 10000046
$ws.Range("C37").Value = 8000046
$ws.Range("D37").Value = "IRPI CNR"
$ws.Range("E37").Value = "Modelli e carte di suscettibilità da frana"

# Row 38
$ws.Range("B38").Value = 10000050
$ws.Range("D38").Value = "Confederazione elvetica"
$ws.Range("E38").Value = "Legge federale sulla protezione della popolazione e sulla protezione civile del 4 ottobre 2002"
